# The sheet "Hortaliza, Femacal de La Calera - Ajo" gets a new weekly price
# record inserted as row 176 (pushing the existing rows 176-257 down to
# 177-258, i.e. dimension grows from A1:R257 to A1:R258).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 176; this shifts rows 176:257 down to 177:258
# and carries the column D date-format styling onto the new row, same as
# native Excel "Insert Row" behaviour.
$ws.Rows("176:176").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(176, 1).Value = 3
$ws.Cells.Item(176, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(176, 3).Value = "Coquimbo"
$ws.Cells.Item(176, 4).Value = 44523
$ws.Cells.Item(176, 5).Value = 5
$ws.Cells.Item(176, 6).Value = 100112003
$ws.Cells.Item(176, 7).Value = "Ajo"
$ws.Cells.Item(176, 8).Value = "Chino"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 73
$ws.Cells.Item(176, 11).Value = 16000
$ws.Cells.Item(176, 12).Value = 16500
$ws.Cells.Item(176, 13).Value = 16240
$ws.Cells.Item(176, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(176, 15).Value = "China"
$ws.Cells.Item(176, 16).Value = 1624
$ws.Cells.Item(176, 17).Value = 10
$ws.Cells.Item(176, 18).Value = "Hortaliza"
